$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for the new "Area" / "Atotal" columns
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"

# G2: first area segment has no preceding depth, so anchor at 0
$ws.Range("G2").Formula = "=(D2-0)*B2/100"

# G3: first real incremental area
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"

# G4:G15 share the general incremental-area formula (fill down as one shared formula)
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# Total area
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# Match the author's cursor position after adding the new columns
$ws.Range("H2").Select()

$wb.Save()
